$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 736
$sheet1.Range("F3").Value = 68
$sheet1.Range("F5").Value = 6069
$sheet1.Range("F7").Value = 190
$sheet1.Range("F12").Value = 4725
$sheet1.Range("F16").Value = 70
$sheet1.Range("F19").Value = 1054
$sheet1.Range("F23").Value = 93
$sheet1.Range("F27").Value = 9
$sheet1.Range("F30").Value = 123
$sheet1.Range("F39").Value = 298
$sheet1.Range("F40").Value = 45
$sheet1.Range("F43").Value = 59
$sheet1.Range("F46").Value = 468
$sheet1.Range("F47").Value = 472

$sheet2 = $wb.Worksheets.Item("演出")
$sheet2.Range("F21").Value = 1
$sheet2.Range("F23").Value = 7

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 736
$sheet4.Range("F3").Value = 68
$sheet4.Range("F5").Value = 6069
$sheet4.Range("F7").Value = 190
$sheet4.Range("F11").Value = 4725
$sheet4.Range("F19").Value = 70
$sheet4.Range("F23").Value = 1054
$sheet4.Range("F32").Value = 123
$sheet4.Range("F43").Value = 298
$sheet4.Range("F44").Value = 1
$sheet4.Range("F47").Value = 468
$sheet4.Range("F48").Value = 472
$sheet4.Range("F49").Value = 7
